# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.335.30'
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = '1.661.88'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.12'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.95'
$ws.Range("E10").Value = '  +4.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0850'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '1.893.88'
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").Value = '1.660.88'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.536'
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.44'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("D17").Value = '27.314.10'
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.84'
$ws.Range("E19").Value = '  +3.54%  '
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.77'
$ws.Range("E21").Value = '  +8.63%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("E23").Value = '  +4.54%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.46'
$ws.Range("E27").Value = '  +3.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.08'
$ws.Range("E29").Value = '  +2.79%  '
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("D35").Value = '1.262.57'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.541'
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.814'
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").Value = '1.805.84'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.13'
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.95'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.44'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0985'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.408'
